$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I19").Value = 1373.8182
$ws.Range("L19").Value = 1138.6666
$ws.Range("M19").Value = -1198.8182
$ws.Range("K19").Value = 1373.8182
$ws.Range("N19").Value = -1488.6666
$ws.Range("J19").Value = 1138.6666
$ws.Range("H19").Value = 1268
$ws.Range("H62").Value = 7094.2144
$ws.Range("I62").Value = 4914.125
$ws.Range("M62").Value = -4290.125
$ws.Range("K62").Value = 4914.125
$ws.Range("K65").Value = 24570.625
$ws.Range("H65").Value = 7094.2144
$ws.Range("I65").Value = 4914.125
$ws.Range("M65").Value = -21450.625
$ws.Range("N76").Value = -9629.454
$ws.Range("J76").Value = 8999.454
$ws.Range("H76").Value = 8649.9375
$ws.Range("L76").Value = 8999.454
$ws.Range("H79").Value = 8649.9375
$ws.Range("L79").Value = 8999.454
$ws.Range("N79").Value = -11183.454
$ws.Range("J79").Value = 8999.454
$ws.Range("I86").Value = 4263
$ws.Range("L86").Value = 4990.8184
$ws.Range("J86").Value = 4990.8184
$ws.Range("M86").Value = -3140
$ws.Range("H86").Value = 4684.3687
$ws.Range("K86").Value = 4263
$ws.Range("N86").Value = -7236.8184
$ws.Range("N89").Value = -36186.092
$ws.Range("J89").Value = 4990.8184
$ws.Range("H89").Value = 4684.3687
$ws.Range("I89").Value = 4263
$ws.Range("L89").Value = 24954.092
$ws.Range("M89").Value = -15699
$ws.Range("K89").Value = 21315
$ws.Range("M92").Value = -1068.611
$ws.Range("K92").Value = 2316.611
$ws.Range("H92").Value = 2963.0952
$ws.Range("I92").Value = 2316.611
$ws.Range("H94").Value = 3436.75
$ws.Range("I94").Value = 3356.5715
$ws.Range("M94").Value = -2905.5715
$ws.Range("K94").Value = 3356.5715
$ws.Range("H107").Value = 4629.4
$ws.Range("I107").Value = 4548.8335
$ws.Range("M107").Value = -2628.8335
$ws.Range("K107").Value = 4548.8335
$ws.Range("I116").Value = 3995
$ws.Range("L116").Value = 12476.5
$ws.Range("J116").Value = 12476.5
$ws.Range("M116").Value = -553
$ws.Range("H116").Value = 10780.2
$ws.Range("K116").Value = 3995
$ws.Range("N116").Value = -19360.5
$ws.Range("I132").Value = 1538.3429
$ws.Range("M132").Value = -2085.028700000001
$ws.Range("H132").Value = 2419.1025
$ws.Range("K132").Value = 4615.028700000001
$ws.Range("N135").Value = -20083.4994
$ws.Range("J135").Value = 1668.1666
$ws.Range("H135").Value = 1753.9
$ws.Range("L135").Value = 15013.4994
$ws.Range("H137").Value = 3155.1912
$ws.Range("I137").Value = 2536.75
$ws.Range("M137").Value = -5060.25
$ws.Range("K137").Value = 7610.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1324.4728
$ws.Range("I32").Value = 1403.1224
$ws.Range("M32").Value = -1116.1224
$ws.Range("K32").Value = 1403.1224
$ws.Range("K61").Value = 5305.6875
$ws.Range("H61").Value = 7351.952
$ws.Range("I61").Value = 5305.6875
$ws.Range("M61").Value = -5093.6875
$ws.Range("J74").Value = 4749.8335
$ws.Range("L74").Value = 4749.8335
$ws.Range("H74").Value = 12824920
$ws.Range("N74").Value = -6497.8335
$ws.Range("J77").Value = 4749.8335
$ws.Range("H77").Value = 12824920
$ws.Range("L77").Value = 23749.1675
$ws.Range("N77").Value = -32485.1675
$ws.Range("K97").Value = 6725.231
$ws.Range("N97").Value = -9930.25
$ws.Range("I97").Value = 6725.231
$ws.Range("J97").Value = 8938.25
$ws.Range("M97").Value = -6229.231
$ws.Range("H97").Value = 7245.9414
$ws.Range("L97").Value = 8938.25
$ws.Range("L102").Value = 2600
$ws.Range("J102").Value = 2600
$ws.Range("H102").Value = 2664.2856
$ws.Range("I102").Value = 2690
$ws.Range("M102").Value = -1068
$ws.Range("K102").Value = 2690
$ws.Range("N102").Value = -5844
$ws.Range("K122").Value = 10269.3
$ws.Range("H122").Value = 3738.48
$ws.Range("I122").Value = 3423.1
$ws.Range("M122").Value = -7819.299999999999
$ws.Range("I132").Value = 5346.846
$ws.Range("M132").Value = -13510.538
$ws.Range("H132").Value = 6453.3228
$ws.Range("K132").Value = 16040.538
$ws.Range("I136").Value = 5305.6875
$ws.Range("M136").Value = -13367.0625
$ws.Range("H136").Value = 7351.952
$ws.Range("K136").Value = 15917.0625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I20").Value = 1688.25
$ws.Range("M20").Value = -1441.25
$ws.Range("H20").Value = 1507.3077
$ws.Range("K20").Value = 1688.25
$ws.Range("I86").Value = 1598.7142
$ws.Range("L86").Value = 5901.75
$ws.Range("J86").Value = 5901.75
$ws.Range("M86").Value = -475.7141999999999
$ws.Range("H86").Value = 3163.4546
$ws.Range("K86").Value = 1598.7142
$ws.Range("N86").Value = -8147.75
$ws.Range("J88").Value = 26349.375
$ws.Range("H88").Value = 26349.375
$ws.Range("L88").Value = 26349.375
$ws.Range("N88").Value = -27161.375
$ws.Range("N89").Value = -40740.75
$ws.Range("J89").Value = 5901.75
$ws.Range("H89").Value = 3163.4546
$ws.Range("I89").Value = 1598.7142
$ws.Range("L89").Value = 29508.75
$ws.Range("M89").Value = -2377.571
$ws.Range("K89").Value = 7993.571
$ws.Range("N91").Value = -29157.375
$ws.Range("J91").Value = 26349.375
$ws.Range("H91").Value = 26349.375
$ws.Range("L91").Value = 26349.375
$ws.Range("H94").Value = 840.4
$ws.Range("I94").Value = 791.4545000000001
$ws.Range("M94").Value = -340.4545000000001
$ws.Range("K94").Value = 791.4545000000001
$ws.Range("K99").Value = 5000.5
$ws.Range("H99").Value = 5000.5
$ws.Range("I99").Value = 5000.5
$ws.Range("M99").Value = -3502.5
$ws.Range("J105").Value = 11047.4
$ws.Range("H105").Value = 16530.117
$ws.Range("I105").Value = 18814.584
$ws.Range("L105").Value = 11047.4
$ws.Range("M105").Value = -17067.584
$ws.Range("K105").Value = 18814.584
$ws.Range("N105").Value = -14541.4
$ws.Range("J134").Value = 10342.333
$ws.Range("H134").Value = 2598.5557
$ws.Range("I134").Value = 1630.5834
$ws.Range("L134").Value = 31026.999
$ws.Range("M134").Value = -2356.7502
$ws.Range("K134").Value = 4891.7502
$ws.Range("N134").Value = -36096.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 430298.72
$ws.Range("I6").Value = 430298.72
$ws.Range("M6").Value = -430185.72
$ws.Range("K6").Value = 430298.72
$ws.Range("N31").Value = -33889
$ws.Range("K31").Value = 3161
$ws.Range("J31").Value = 33299
$ws.Range("H31").Value = 25025.824
$ws.Range("I31").Value = 3161
$ws.Range("L31").Value = 33299
$ws.Range("M31").Value = -2866
$ws.Range("J34").Value = 33299
$ws.Range("H34").Value = 25025.824
$ws.Range("I34").Value = 3161
$ws.Range("L34").Value = 33299
$ws.Range("M34").Value = -2959
$ws.Range("K34").Value = 3161
$ws.Range("N34").Value = -33703
$ws.Range("L43").Value = 10000
$ws.Range("J43").Value = 10000
$ws.Range("H43").Value = 10000
$ws.Range("N43").Value = -10368
$ws.Range("L58").Value = 21008
$ws.Range("J58").Value = 21008
$ws.Range("H58").Value = 4966.95
$ws.Range("N58").Value = -21414
$ws.Range("I86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H86").Value = 15817.333
$ws.Range("K86").Value = 0
$ws.Range("H89").Value = 15817.333
$ws.Range("I89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("K89").Value = 0
$ws.Range("J101").Value = 10000
$ws.Range("L101").Value = 10000
$ws.Range("H101").Value = 10000
$ws.Range("N101").Value = -16490
$ws.Range("H103").Value = 524
$ws.Range("I103").Value = 524
$ws.Range("M103").Value = 648
$ws.Range("K103").Value = 524
$ws.Range("H105").Value = 2765.1365
$ws.Range("I105").Value = 3147
$ws.Range("M105").Value = -1400
$ws.Range("K105").Value = 3147
$ws.Range("K122").Value = 6456.882599999999
$ws.Range("N122").Value = -55540
$ws.Range("J122").Value = 16880
$ws.Range("H122").Value = 5499.5
$ws.Range("I122").Value = 2152.2942
$ws.Range("L122").Value = 50640
$ws.Range("M122").Value = -4006.882599999999
$ws.Range("H134").Value = 3018.1785
$ws.Range("I134").Value = 1770.625
$ws.Range("M134").Value = -2776.875
$ws.Range("K134").Value = 5311.875
$ws.Range("N136").Value = -68124
$ws.Range("J136").Value = 21008
$ws.Range("L136").Value = 63024
$ws.Range("H136").Value = 4966.95

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 19309.8
$ws.Range("I32").Value = 1700
$ws.Range("M32").Value = -4817
$ws.Range("K32").Value = 5100
$ws.Range("J34").Value = 2224.625
$ws.Range("H34").Value = 1497.6428
$ws.Range("I34").Value = 528.3333
$ws.Range("L34").Value = 6673.875
$ws.Range("M34").Value = -1500.9999
$ws.Range("K34").Value = 1584.9999
$ws.Range("N34").Value = -6841.875
$ws.Range("H47").Value = 174998.5
$ws.Range("I47").Value = 100000
$ws.Range("M47").Value = -299569
$ws.Range("K47").Value = 300000
$ws.Range("I58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("H87").Value = 12890.223
$ws.Range("I87").Value = 10602.4
$ws.Range("M87").Value = -30559.2
$ws.Range("K87").Value = 31807.2
$ws.Range("I90").Value = 10602.4
$ws.Range("M90").Value = -89181.59999999999
$ws.Range("H90").Value = 12890.223
$ws.Range("K90").Value = 95421.59999999999
$ws.Range("J113").Value = 1607.35
$ws.Range("L113").Value = 4822.049999999999
$ws.Range("H113").Value = 1432.3793
$ws.Range("I113").Value = 1043.5555
$ws.Range("M113").Value = -960.6664999999998
$ws.Range("K113").Value = 3130.6665
$ws.Range("N113").Value = -9162.049999999999
$ws.Range("L114").Value = 5383.5
$ws.Range("M114").Value = -2546.0002
$ws.Range("K114").Value = 5800.0002
$ws.Range("N114").Value = -11891.5
$ws.Range("J114").Value = 1794.5
$ws.Range("H114").Value = 1854
$ws.Range("I114").Value = 1933.3334
$ws.Range("J117").Value = 4501.125
$ws.Range("H117").Value = 3717
$ws.Range("I117").Value = 2671.5
$ws.Range("L117").Value = 13503.375
$ws.Range("M117").Value = -4572.5
$ws.Range("K117").Value = 8014.5
$ws.Range("N117").Value = -20387.375
$ws.Range("J121").Value = 1409.3636
$ws.Range("L121").Value = 4228.0908
$ws.Range("H121").Value = 1807.9231
$ws.Range("N121").Value = -6848.0908
$ws.Range("J129").Value = 7579925.5
$ws.Range("H129").Value = 5558924
$ws.Range("L129").Value = 22739776.5
$ws.Range("N129").Value = -22749776.5
$ws.Range("I132").Value = 4017.5715
$ws.Range("J132").Value = 6592.3335
$ws.Range("L132").Value = 59331.0015
$ws.Range("M132").Value = -33628.1435
$ws.Range("H132").Value = 5205.923
$ws.Range("K132").Value = 36158.1435
$ws.Range("N132").Value = -64391.0015
$ws.Range("M139").Value = -1660.000100000001
$ws.Range("H139").Value = 6395.45
$ws.Range("I139").Value = 2266.6667
$ws.Range("K139").Value = 6800.000100000001
$ws.Range("H140").Value = 6004.25
$ws.Range("I140").Value = 4495
$ws.Range("M140").Value = -8305
$ws.Range("K140").Value = 13485

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8411.5
$ws.Range("I70").Value = 7458.8
$ws.Range("M70").Value = -7188.8
$ws.Range("K70").Value = 7458.8
$ws.Range("K73").Value = 7458.8
$ws.Range("H73").Value = 8411.5
$ws.Range("I73").Value = 7458.8
$ws.Range("M73").Value = -6522.8
$ws.Range("L80").Value = 15425.25
$ws.Range("M80").Value = -3885
$ws.Range("N80").Value = -17421.25
$ws.Range("K80").Value = 4883
$ws.Range("J80").Value = 15425.25
$ws.Range("H80").Value = 9099.9
$ws.Range("I80").Value = 4883
$ws.Range("I83").Value = 4883
$ws.Range("J83").Value = 15425.25
$ws.Range("L83").Value = 77126.25
$ws.Range("M83").Value = -19423
$ws.Range("H83").Value = 9099.9
$ws.Range("K83").Value = 24415
$ws.Range("N83").Value = -87110.25
$ws.Range("H102").Value = 4474.25
$ws.Range("I102").Value = 2900
$ws.Range("M102").Value = -1278
$ws.Range("K102").Value = 2900
$ws.Range("H107").Value = 1758.579
$ws.Range("I107").Value = 701.6
$ws.Range("L107").Value = 5722.25
$ws.Range("M107").Value = 1218.4
$ws.Range("K107").Value = 701.6
$ws.Range("N107").Value = -9562.25
$ws.Range("J107").Value = 5722.25
$ws.Range("K122").Value = 23845.6158
$ws.Range("H122").Value = 8361.143
$ws.Range("I122").Value = 7948.5386
$ws.Range("M122").Value = -21395.6158
$ws.Range("J126").Value = 9228.1
$ws.Range("H126").Value = 6014.778
$ws.Range("L126").Value = 27684.3
$ws.Range("N126").Value = -32624.3
$ws.Range("I132").Value = 3881.9333
$ws.Range("M132").Value = -9115.7999
$ws.Range("H132").Value = 5250.4116
$ws.Range("K132").Value = 11645.7999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K40").Value = 7850.9414
$ws.Range("N40").Value = -8971.200000000001
$ws.Range("I40").Value = 7850.9414
$ws.Range("J40").Value = 8699.200000000001
$ws.Range("M40").Value = -7714.9414
$ws.Range("L40").Value = 8699.200000000001
$ws.Range("H40").Value = 8165.1113
$ws.Range("K82").Value = 4685.8887
$ws.Range("N82").Value = -7721.875
$ws.Range("J82").Value = 6999.875
$ws.Range("H82").Value = 5774.8237
$ws.Range("I82").Value = 4685.8887
$ws.Range("L82").Value = 6999.875
$ws.Range("M82").Value = -4324.8887
$ws.Range("N85").Value = -9495.875
$ws.Range("L85").Value = 6999.875
$ws.Range("J85").Value = 6999.875
$ws.Range("H85").Value = 5774.8237
$ws.Range("I85").Value = 4685.8887
$ws.Range("M85").Value = -3437.8887
$ws.Range("K85").Value = 4685.8887
$ws.Range("K93").Value = 11695.381
$ws.Range("N93").Value = -20515.928
$ws.Range("I93").Value = 11695.381
$ws.Range("J93").Value = 18019.928
$ws.Range("M93").Value = -10447.381
$ws.Range("L93").Value = 18019.928
$ws.Range("H93").Value = 14225.2
$ws.Range("K100").Value = 851
$ws.Range("H100").Value = 6902
$ws.Range("I100").Value = 851
$ws.Range("M100").Value = -310
$ws.Range("K122").Value = 65457.702
$ws.Range("H122").Value = 21470.105
$ws.Range("I122").Value = 21819.234
$ws.Range("M122").Value = -63007.702
$ws.Range("I132").Value = 3105.1875
$ws.Range("J132").Value = 13252.5
$ws.Range("L132").Value = 39757.5
$ws.Range("M132").Value = -6785.5625
$ws.Range("H132").Value = 5134.65
$ws.Range("K132").Value = 9315.5625
$ws.Range("N132").Value = -44817.5
$ws.Range("N136").Value = -43982.7
$ws.Range("I136").Value = 5026.923
$ws.Range("J136").Value = 12960.9
$ws.Range("M136").Value = -12530.769
$ws.Range("L136").Value = 38882.7
$ws.Range("H136").Value = 8476.478999999999
$ws.Range("K136").Value = 15080.769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 50000
$ws.Range("N46").Value = -50462
$ws.Range("J46").Value = 50000
$ws.Range("L46").Value = 50000
$ws.Range("H81").Value = 7922.923
$ws.Range("I81").Value = 5077.1113
$ws.Range("L81").Value = 28652
$ws.Range("M81").Value = -9093.222599999999
$ws.Range("N81").Value = -30774
$ws.Range("K81").Value = 10154.2226
$ws.Range("J81").Value = 14326
$ws.Range("H84").Value = 7922.923
$ws.Range("K84").Value = 50771.113
$ws.Range("N84").Value = -153868
$ws.Range("J84").Value = 14326
$ws.Range("I84").Value = 5077.1113
$ws.Range("L84").Value = 143260
$ws.Range("M84").Value = -45467.113
$ws.Range("K100").Value = 1743.5
$ws.Range("H100").Value = 1181.6666
$ws.Range("I100").Value = 871.75
$ws.Range("M100").Value = -1202.5
$ws.Range("J113").Value = 1362.3334
$ws.Range("L113").Value = 4087.0002
$ws.Range("H113").Value = 1220
$ws.Range("I113").Value = 1107.6316
$ws.Range("M113").Value = -1152.8948
$ws.Range("K113").Value = 3322.8948
$ws.Range("N113").Value = -8427.0002
$ws.Range("K122").Value = 6171.428400000001
$ws.Range("N122").Value = -31539.4
$ws.Range("J122").Value = 8879.799999999999
$ws.Range("H122").Value = 3852.5789
$ws.Range("I122").Value = 2057.1428
$ws.Range("L122").Value = 26639.4
$ws.Range("M122").Value = -3721.428400000001
$ws.Range("I132").Value = 3270.3794
$ws.Range("M132").Value = -7281.138199999999
$ws.Range("H132").Value = 6407.2354
$ws.Range("K132").Value = 9811.138199999999
$ws.Range("J134").Value = 50000
$ws.Range("H134").Value = 50000
$ws.Range("L134").Value = 150000
$ws.Range("N134").Value = -155070
$ws.Range("N136").Value = -28672.7139
$ws.Range("I136").Value = 5724.75
$ws.Range("J136").Value = 7857.5713
$ws.Range("M136").Value = -14624.25
$ws.Range("L136").Value = 23572.7139
$ws.Range("H136").Value = 7082
$ws.Range("K136").Value = 17174.25
